# Update the existing "Step 2" cell text and its row's expected-result text,
# then add the new "Step 3" row (Steps + Expected Result columns).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Step 2 text is updated, and the expected result for it is now
# "I am redirected to the user's dashboard".
$ws.Range("C3").Value = "Step 2: Log in as a user with the appropriate role"
$ws.Range("D3").Value = "I am redirected to the user's dashboard"

# Row 4: new Step 3 with its expected result.
$ws.Range("C4").Value = "Step 3: Click on the ""Employees"" Button"
$ws.Range("D4").Value = "I see a list of employees of onshore"
